# Add a new "FS" (factor de seguridad) column to the foundation data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in column I, row 1, and its value in row 2.
$ws.Range("I1").Value = "FS"
$ws.Range("I2").Value = 3

# Match the formatting already used by the rest of the header/data row
# (centered horizontal + vertical alignment) by copying it from the
# neighbouring column H instead of building a brand-new style.
$ws.Range("H1:H2").Copy()
$ws.Range("I1:I2").PasteSpecial(-4122)  # xlPasteFormats

# Move the active selection to I3, just below the new column, as in the
# original author's edit.
[void]$ws.Range("I3").Select()
